$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Change alignment of rows 7 and 20 from Left to Right
$ws.Range("D7").Value = "Right"
$ws.Range("D20").Value = "Right"

# Add new interaction rows 41-44
$ws.Range("B41").Value = "SingleUseId64"
$ws.Range("C41").Value = "Small"
$ws.Range("D41").Value = "Center"
$ws.Range("E41").Value = "LTR"
$ws.Range("F41").Value = "Yes"

$ws.Range("B42").Value = "SingleUseId65"
$ws.Range("C42").Value = "Small"
$ws.Range("D42").Value = "Center"
$ws.Range("E42").Value = "LTR"
$ws.Range("F42").Value = "Cancel"

$ws.Range("B43").Value = "SingleUseId66"
$ws.Range("C43").Value = "Default"
$ws.Range("D43").Value = "Left"
$ws.Range("E43").Value = "LTR"
$ws.Range("F43").Value = "Reset timer?"

$ws.Range("B44").Value = "SingleUseId67"
$ws.Range("C44").Value = "Tiny"
$ws.Range("D44").Value = "Right"
$ws.Range("E44").Value = "LTR"
$ws.Range("F44").Value = "Limit"
